$wb = $excel.ActiveWorkbook

# --- Sheet1: append new TODO items in column A, rows 10-14 ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$newItems = @(
    "RefSource selection doesný need to be on output page, put it to Source page as checkbox or something (only needed for RefSync mode - no need to be on outpt page that is for all)",
    "Make Source name edit focused control when new Source is created",
    "Custom naming patterns - as advanced option",
    "Add seconds to default pattern - when adding additional prhotos to previously done mix there maybe overwrites since the counter starts at zero on additional mix",
    "When Tool is executed go to output page to see the log"
)

$row = 10
foreach ($item in $newItems) {
    $ws1.Cells.Item($row, 1).Value = $item
    $row++
}

# Update active selection on Sheet1 to A14 (last added row)
$ws1.Activate()
$ws1.Range("A14").Select()

# --- Sheet2: set active selection to B27 ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("B27").Select()

# Re-activate Sheet1 as the originally selected/tabbed sheet
$ws1.Activate()
